$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: new "APP" row, bold red font, bordered (matches existing table border style) ---
# Seed the new cells with the existing bordered style (same as rest of the table) by
# copying format from an already-bordered, unfilled cell, then apply bold red font.
$ws.Range("B7").Copy() | Out-Null
$ws.Range("A9:C9").PasteSpecial(-4122) | Out-Null
$ws.Range("E9:F9").PasteSpecial(-4122) | Out-Null

$ws.Range("A9").Value = "Cplayer APP"
$ws.Range("B9").Value = "MUZIX APP"
$ws.Range("C9").Value = "WALMART APP"
$ws.Range("F9").Value = "NEWS APP"
$ws.Range("E9").Value = "TRANSPORT APP"

$ws.Range("A9:C9").Font.Bold = $true
$ws.Range("A9:C9").Font.Color = 255
$ws.Range("E9:F9").Font.Bold = $true
$ws.Range("E9:F9").Font.Color = 255

# --- Rows 12-15: small design/testing notes block ---
$ws.Range("B7").Copy() | Out-Null
$ws.Range("D12:E15").PasteSpecial(-4122) | Out-Null

$ws.Range("D12").Value = "1)Design "
$ws.Range("D13").Value = "2)Document"
$ws.Range("E12").Value = "usecase,DFD,class,forms,tables"
$ws.Range("E13").Value = "Modification in the abstract,s/w requirements"
$ws.Range("D14").Value = "3)coding"
$ws.Range("D15").Value = "4)testing"
$ws.Range("E15").Value = "junit testing"

$ws.Range("E13").WrapText = $true
$ws.Rows(13).RowHeight = 30

# D9 gets the same bold red font plus an explicit (no-op) fill toggle; entered last to
# mirror the original authoring order (WEATHER APP was the final string added).
$ws.Range("B7").Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4122) | Out-Null
$ws.Range("D9").Interior.Pattern = -4142
$ws.Range("D9").Interior.PatternColorIndex = -4105
$ws.Range("D9").Value = "WEATHER APP"
$ws.Range("D9").Font.Bold = $true
$ws.Range("D9").Font.Color = 255

# --- Column / selection / dimension cosmetics ---
$ws.Columns("E").ColumnWidth = 29.5
$ws.Range("F9").Select() | Out-Null
